$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the date cell's format (style index) onto the new date cells first
$ws.Range("A2").Copy()
$ws.Range("A17:A18").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("A17").Value = 42602.010115740741
$ws.Range("B17").Value = "Noun"
$ws.Range("C17").Value = 7613
$ws.Range("D17").Value = 4193
$ws.Range("E17").Value = 815
$ws.Range("F17").Value = 130
$ws.Range("G17").Value = 58
$ws.Range("H17").Value = 68
$ws.Range("I17").Value = 30
$ws.Range("J17").Value = 0
$ws.Range("K17").Value = 0
$ws.Range("L17").Value = 0
$ws.Range("M17").Value = 0

$ws.Range("A18").Value = 42602.481909722221
$ws.Range("B18").Value = "Noun"
$ws.Range("C18").Value = 8856
$ws.Range("D18").Value = 6064
$ws.Range("E18").Value = 1179
$ws.Range("F18").Value = 193
$ws.Range("G18").Value = 100
$ws.Range("H18").Value = 65
$ws.Range("I18").Value = 34
$ws.Range("J18").Value = 1
$ws.Range("K18").Value = 0
$ws.Range("L18").Value = 100
$ws.Range("M18").Value = 0
